$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A width change (14 -> 10) ---
$ws.Columns.Item(1).ColumnWidth = 9.17

# --- Line item rows (2-4): new partidas / quantities / rates / time / subtotal ---
$ws.Range("A2").Value = "Igenieros"
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 2000
$ws.Range("D2").Value = 2
$ws.Range("E2").Value = 40000

$ws.Range("A3").Value = "aaa"
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 1200
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 4800

$ws.Range("A4").Value = "xd"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 500
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 1000

# --- Old row 5 (Estimacion 4) and row 6 (Estimacion 5) data items are gone ---
$ws.Range("A5:E6").ClearContents()

# --- Bring the shaded/bordered "totals column" style (used by E8) onto the
#     new summary rows E6/E7, which previously held plain unstyled values ---
$ws.Range("E8").Copy($ws.Range("E6"))
$ws.Range("E8").Copy($ws.Range("E7"))

# --- Summary block now starts at row 6 instead of row 8 ---
$ws.Range("D6").Value = "TOTAL"
$ws.Range("E6").Value = 45800

$ws.Range("D7").Value = "Reserva de contingencia"
$ws.Range("E7").Value = 1500

$ws.Range("D8").Value = "Linea Base de Costos"
$ws.Range("E8").Value = 47300

$ws.Range("D9").Value = "Reserva de gestion"
$ws.Range("E9").Value = "PV"
$ws.Range("F9").Value = 0.2

$ws.Range("D10").Value = "Presupuesto"
$ws.Range("E10").Value = 12500

$ws.Range("D11").Value = "Ganancia"
$ws.Range("E11").Value = "PV"
$ws.Range("F11").Value = 0.4

$ws.Range("D12").Value = "Total con ganancia"
$ws.Range("E12").Value = 59800

$ws.Range("D13").Value = "IGV"
$ws.Range("E13").Value = "PV"
$ws.Range("F13").Value = 0.5

$ws.Range("D14").Value = "Total"
$ws.Range("E14").Value = 59800.5

# --- Former rows 15 and 16 no longer exist; remove them outright so the
#     sheet dimension shrinks back to row 14 ---
$ws.Range("A15:F16").EntireRow.Delete()
